$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (e.g. "1.00", "0.999").
# Excel would otherwise auto-convert these to numbers and lose the exact text
# (trailing zeros, etc.), so temporarily force a Text number format while
# assigning the value, then clear the format again so the cell ends up with
# its original (default) style, just like every other text cell on the sheet.
$textCells = @("D5", "D6", "D8", "D11", "D12", "D13", "D14", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D35", "D38", "D39", "D41", "D42", "D44", "D45", "D46", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '62.696.81'
$ws.Range('E2').Value = '  -4.38%  '
$ws.Range('D3').Value = '3.264.31'
$ws.Range('E3').Value = '  -6.86%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '536.30'
$ws.Range('E5').Value = '  -3.67%  '
$ws.Range('D6').Value = '168.85'
$ws.Range('E6').Value = '  -5.99%  '
$ws.Range('E7').Value = '  -5.45%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = '3.252.83'
$ws.Range('E9').Value = '  -7.04%  '
$ws.Range('E10').Value = '  -4.65%  '
$ws.Range('D11').Value = '0.150'
$ws.Range('E11').Value = '  -3.12%  '
$ws.Range('D12').Value = '51.83'
$ws.Range('E12').Value = '  -3.92%  '
$ws.Range('D13').Value = '0.0000260'
$ws.Range('E13').Value = '  -5.08%  '
$ws.Range('D14').Value = '8.74'
$ws.Range('E14').Value = '  -5.75%  '
$ws.Range('D15').Value = '3.798.04'
$ws.Range('E15').Value = '  -6.59%  '
$ws.Range('E16').Value = '  -3.50%  '
$ws.Range('E17').Value = '  -4.58%  '
$ws.Range('D18').Value = '3.264.15'
$ws.Range('E18').Value = '  -6.88%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '11.47'
$ws.Range('E19').Value = '  -5.61%  '
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = '62.591.72'
$ws.Range('E20').Value = '  -4.52%  '
$ws.Range('D21').Value = '0.955'
$ws.Range('E21').Value = '  -4.31%  '
$ws.Range('D22').Value = '408.99'
$ws.Range('E22').Value = '  -1.81%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').Value = '4.32'
$ws.Range('E23').Value = '  +4.83%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Value = '3.95'
$ws.Range('E24').Value = '  -2.57%  '
$ws.Range('D25').Value = '13.29'
$ws.Range('E25').Value = '  +3.29%  '
$ws.Range('D26').Value = '82.15'
$ws.Range('E26').Value = '  -4.63%  '
$ws.Range('D27').Value = '10.43'
$ws.Range('E27').Value = '  -3.69%  '
$ws.Range('E28').Value = '  -6.30%  '
$ws.Range('D29').Value = '8.44'
$ws.Range('E29').Value = '  -6.98%  '
$ws.Range('D30').Value = '28.67'
$ws.Range('E30').Value = '  -5.72%  '
$ws.Range('D31').Value = '6.23'
$ws.Range('E31').Value = '  -4.14%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = '11.18'
$ws.Range('E32').Value = '  -4.41%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = '566.97'
$ws.Range('E33').Value = '  -6.74%  '
$ws.Range('E34').Value = '  -4.91%  '
$ws.Range('D35').Value = '57.58'
$ws.Range('E35').Value = '  -3.63%  '
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('E37').Value = '  -2.35%  '
$ws.Range('D38').Value = '34.60'
$ws.Range('E38').Value = '  -7.47%  '
$ws.Range('D39').Value = '3.36'
$ws.Range('E39').Value = '  +3.25%  '
$ws.Range('D40').Value = '0.0₃0725'
$ws.Range('E40').Value = '  -8.89%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').Value = '0.358'
$ws.Range('E41').Value = '  -6.00%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('D43').Value = '3.072.46'
$ws.Range('E43').Value = '  -9.26%  '
$ws.Range('D44').Value = '3.21'
$ws.Range('E44').Value = '  -1.01%  '
$ws.Range('D45').Value = '2.69'
$ws.Range('E45').Value = '  -5.37%  '
$ws.Range('D46').Value = '0.0393'
$ws.Range('E46').Value = '  -5.30%  '
$ws.Range('E47').Value = '  -6.62%  '
$ws.Range('E49').Value = '  -4.73%  '
$ws.Range('D50').Value = '131.63'
$ws.Range('E50').Value = '  -4.60%  '
$ws.Range('E51').Value = '  -6.65%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
